$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English and Communication", 1634, "Female", 2000, "Intermediate_2"),
    @("French", 461, "Female", 2000, "Intermediate_2"),
    @("Mathematics", 4402, "Female", 2000, "Intermediate_2"),
    @("Biology", 1767, "Female", 2000, "Intermediate_2"),
    @("Chemistry", 630, "Female", 2000, "Intermediate_2"),
    @("Physics", 451, "Female", 2000, "Intermediate_2"),
    @("Computing", 229, "Female", 2000, "Intermediate_2")
)

$row = 9
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
